# Update the "Price" column (D) on the active sheet with refreshed
# cryptocurrency quotes, as produced by the scheduled GitHub Actions
# symbol-list refresh job.
#
# The price values are stored as plain text (not numbers) in the sheet,
# so each cell's number format is forced to Text ("@") before the new
# value is written. This prevents Excel from re-interpreting strings
# such as "0.001658" or "0.00000000752" as numbers (which would mangle
# their literal representation, e.g. via scientific notation or
# trailing-zero loss).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "246.62"
    "D3"  = "22.41"
    "D4"  = "5.506"
    "D5"  = "0.05631"
    "D6"  = "6.469"
    "D7"  = "0.8058"
    "D8"  = "1.053"
    "D9"  = "0.1441"
    "D10" = "0.07345"
    "D11" = "0.03183"
    "D12" = "0.02926"
    "D13" = "0.09248"
    "D14" = "0.001658"
    "D15" = "3.199"
    "D16" = "0.04735"
    "D17" = "0.0005847"
    "D18" = "0.006304"
    "D19" = "0.001056"
    "D21" = "0.0001504"
    "D22" = "3.981"
    "D23" = "3.384"
    "D24" = "2.133"
    "D25" = "0.3275"
    "D26" = "0.1316"
    "D27" = "0.0003008"
    "D40" = "0.04159"
    "D41" = "0.006877"
    "D42" = "0.1038"
    "D43" = "0.003202"
    "D44" = "0.009063"
    "D45" = "0.00005662"
    "D46" = "0.00000000752"
    "D47" = "0.6818"
    "D48" = "0.02083"
    "D49" = "0.00002105"
    "D50" = "0.01013"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
